$d = $word.ActiveDocument

# Italicize every "{{TIME}}" placeholder in the document body.
# Word splits the run(s) it touches automatically, so the run that used
# to hold "<tab/>{{TIME}}" becomes a tab run followed by a new,
# separately-formatted "{{TIME}}" run - matching the target diff, which
# turns:
#   <w:r><w:tab/><w:t>{{TIME}}</w:t></w:r>
# into:
#   <w:r><w:tab/></w:r>
#   <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>{{TIME}}</w:t></w:r>

$searchRange = $d.Content
$searchRange.Start = 0

while ($true) {
    $found = $searchRange.Find.Execute("{{TIME}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $searchRange.Font.Italic = $true
    $searchRange.Font.ItalicBi = $true

    # Continue searching right after this match so we pick up every
    # occurrence without re-matching the one we just formatted.
    $nextStart = $searchRange.End
    $searchRange = $d.Content
    $searchRange.Start = $nextStart
}
